# companies_usa_remote.xlsx — refresh scraped rows to the fully
# filtered/deduplicated dataframe: updated elapsed-time + skills values
# for the existing rows, row 10 replaced by the (deduped) Ryder System
# listing, and the trailing duplicate Ryder System rows (11-13) removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header: "Checked companies" count drops 12 -> 9 -----------------
$ws.Range("T1").Value = 9

# --- row 2 : Jobgether / Capacity Manager -----------------------------
$ws.Range("K2").Value = 16.67
$ws.Range("L2").Value = "manage, job, team, process, resource, capacity, match, role, plan, experience"

# --- row 3 : ClearCaptions, LLC / Data Analyst ------------------------
$ws.Range("K3").Value = 20.61
$ws.Range("L3").Value = "data, com, business, work, res, able, team, skill, analysis, unit"

# --- row 4 : Ride Health / Workforce Analyst --------------------------
$ws.Range("K4").Value = 46.57
$ws.Range("L4").Value = "age, com, health, work, ride, workforce, time, manage, management, skill"

# --- row 5 : UNFI / FP & A Analyst II ---------------------------------
$ws.Range("K5").Value = 57.33
$ws.Range("L5").Value = "business, financial, required, remote, experience, work, unfi, team, able, office"

# --- row 6 : TieTalent / Business Metrics/Analytics -------------------
$ws.Range("K6").Value = 65.71
$ws.Range("L6").Value = "management, data, experience, business, work, metrics, portfolio, skills, ability, able"

# --- row 7 : Kforce Inc / Finance Manager, Customer Finance -----------
$ws.Range("K7").Value = 69.11
$ws.Range("L7").Value = "age, finance, service, customer, kforce, financial, team, pay, act, employee"

# --- row 8 : Centene Corporation / Capacity Planning Analyst II -------
$ws.Range("K8").Value = 83.14
$ws.Range("L8").Value = "per, act, capacity, planning, work, center, contact, experience, perform, manage"

# --- row 9 : Jobgether / Workforce Analyst ----------------------------
$ws.Range("K9").Value = 91.71
$ws.Range("L9").Value = "work, per, staffing, job, workforce, team, match, time, teams, manage"

# --- row 10 : was a duplicate Jobgether FP&A posting, now the deduped
#              Ryder System, Inc. / Finance Manager listing. The
#              visa/relocation/remote/applied flags and Job URL already
#              match the Ryder rows so only company, title, elapsed
#              time and skills need to change. -------------------------
$ws.Range("A10").Value = "Ryder System, Inc."
$ws.Range("B10").Value = "Finance Manager - REMOTE"
$ws.Range("K10").Value = 115.76
$ws.Range("L10").Value = "com, age, financial, manage, ryder, app, work, plan, view, job"

# --- rows 11-13 : duplicate Ryder System, Inc. rows dropped by the
#                  dedup/filtered_df logic -----------------------------
$ws.Range("A11:T13").EntireRow.Delete()
